$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.136.81"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -2.62%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.870.52"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -1.91%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.31"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.89%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.02%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5061"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +1.27%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3750"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.64%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07160"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.72%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8902"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.24%  "

$ws.Range("E11").Value = "  -1.22%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.890.15"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.74%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07565"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.00%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.331"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.99%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.45"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.76%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.07%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008514"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -2.51%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.14"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -3.31%  "

$ws.Range("E19").Value = "  +0.15%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.184.63"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.53%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.087"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.75%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.110.45"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.34%  "

$ws.Range("E23").Value = "  -1.84%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.502"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.04%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "150.97"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.46%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.842"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.63%  "

$ws.Range("E27").Value = "  -2.00%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.089"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -6.03%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "112.97"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.04%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.768"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -3.00%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.690"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.68%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08995"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.25%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05137"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.58%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.096"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.97%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7452"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -3.80%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.162"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -5.74%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02038"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.32%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.552"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.34%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.042"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.42%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.078"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.55%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5365"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -3.54%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.618"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -4.15%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "114.96"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.76%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.471"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.26%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1479"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.44%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4645"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.96%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9993"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.01%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.06"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -5.10%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.574"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -3.78%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "64.74"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -4.24%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "36.64"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.33%  "

